$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Metadata" ---
$ws = $wb.Worksheets.Item(1)

# Version bump
$ws.Range("B3").Value = "6.0.0"

# Date update
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher now has a value
$ws.Range("B9").Value = "Alvearie Team"

# The two duplicate "Contact / No display for ContactDetail" rows (10 and 11)
# collapse into a single "Jurisdiction / United States of America" row. Delete
# row 10 so every following row (the rest of the metadata table) shifts up by
# one, keeping all of its original, unmodified content/formatting intact, then
# overwrite the (now blank) row 10 with the new Jurisdiction data.
$ws.Rows.Item(10).Delete()

$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# --- Sheet 2: "Elements" ---
$ws2 = $wb.Worksheets.Item(2)

# Row 2 (the root "Extension" element) gets a proper Short/Definition instead of
# the generic placeholder text
$ws2.Range("K2").Value = "Ethnicity Code"
$ws2.Range("L2").Value = "Code for the ethnicity of the person"
